$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Row 10 (D48 - "Mecanismo para guardar, obtener y mostrar imagenes...") status -> "Hecho"
$ws.Range("F10").Value = "Hecho"

# Row 11 (D50 - "Mostrar unicamente a los alumnos activos para los profesores.") status -> "Hecho"
$ws.Range("F11").Value = "Hecho"

# Register 1 hour consumed on day 5 (column T) for the task in row 11
$ws.Range("T11").Value = 1

# Re-merge the totals header cells so they are re-registered at the end of the
# worksheet's merged-cell collection (matches resulting file ordering)
$ws.Range("AZ4:BA4").UnMerge()
$ws.Range("AZ4:BA4").Merge()
$ws.Range("AO4:AP4").UnMerge()
$ws.Range("AO4:AP4").Merge()
$ws.Range("AR4:AS4").UnMerge()
$ws.Range("AR4:AS4").Merge()
$ws.Range("AU4:AV4").UnMerge()
$ws.Range("AU4:AV4").Merge()
$ws.Range("AX4:AY4").UnMerge()
$ws.Range("AX4:AY4").Merge()

# Update active cell selection to T11 to match the saved view state
$ws.Range("T11").Select()
